# Lesson 0.2 "How to Learn in This Course" - misc typo fixes on slide 12
# (web site / lab notebook paragraph), plus the cached date-placeholder
# text that PowerPoint refreshes ("datetime1"/"datetimeFigureOut" fields)
# on the slide master, every slide layout, and the notes master.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 12 ("Homework policies") - Content Placeholder text edits.
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$contentShape = $null
for ($i = 1; $i -le $s12.Shapes.Count; $i++) {
    $candidate = $s12.Shapes.Item($i)
    if ($candidate.Name -like "Content Placeholder*") {
        $contentShape = $candidate
    }
}

$tr = $contentShape.TextFrame.TextRange

# --- Paragraph 1: split the trailing "." from "... Go study it." into
#     its own run, text itself is unchanged.
$para1Len = 113
$period = $tr.Characters($para1Len, 1)
$period.Text = "."

# --- Paragraph 2: replace the whole "lab notebook" sentence with the
#     new "Work Session Report" text, then split the final two words
#     ("work session.") into a second run.
$len = $tr.Length
$para2Start = $para1Len + 2
$para2Len = $len - $para2Start + 1
$para2 = $tr.Characters($para2Start, $para2Len)
$para2.Text = "The requirement for a lab notebook has been replaced by a Work Session Report, which is a simple Google form that you will complete at the end of each Z"

$newLen = $tr.Length
$splitMarker = $tr.Characters($newLen, 1)
$splitMarker.Text = "work session."

# ---------------------------------------------------------------------
# 2) Refresh the cached "8/1/2015" -> "9/3/2015" date placeholder text
#    on the slide master, every custom layout, and the notes master.
# ---------------------------------------------------------------------
function Update-DateePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = "9/3/2015"
        }
    }
}

$master = $p.SlideMaster
Update-DateePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateePlaceholder $layout.Shapes
}

$notesMaster = $p.NotesMaster
Update-DateePlaceholder $notesMaster.Shapes
